# Updates the "cryptos" price/volume snapshot to the values captured by the
# latest GitHub Actions run. Coin name/link cells are only touched where the
# ranking order of rows 43-48 shifted.
#
# Note: several "Price" values are plain decimal-looking strings (e.g.
# "324.62"). Excel's COM layer auto-detects such text as a number when
# assigned directly, which would change the cell from text to numeric -
# unlike the source data, which always stores these as text. To avoid that,
# those assignments are prefixed with a leading apostrophe, which forces
# Excel to keep (and display) the value as plain text, exactly as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.008.17'
$ws.Range("E2").Value = '  -2.03%  '

$ws.Range("D3").Value = '1.828.89'
$ws.Range("E3").Value = '  -1.12%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = "'324.62"
$ws.Range("E5").Value = '  -2.87%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").Value = "'0.4649"
$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("E8").Value = '  -1.58%  '

$ws.Range("D9").Value = "'0.07852"
$ws.Range("E9").Value = '  -0.73%  '

$ws.Range("D10").Value = "'0.9570"
$ws.Range("E10").Value = '  -2.88%  '

$ws.Range("D11").Value = "'21.84"
$ws.Range("E11").Value = '  -1.77%  '

$ws.Range("D12").Value = '1.887.96'
$ws.Range("E12").Value = '  -3.31%  '

$ws.Range("D13").Value = "'5.671"
$ws.Range("E13").Value = '  -3.06%  '

$ws.Range("D14").Value = "'6.887"
$ws.Range("E14").Value = '  -1.96%  '

$ws.Range("D15").Value = "'0.06852"
$ws.Range("E15").Value = '  -0.18%  '

$ws.Range("D16").Value = "'87.11"
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").Value = "'0.000009905"
$ws.Range("E18").Value = '  -1.64%  '

$ws.Range("D19").Value = "'16.57"
$ws.Range("E19").Value = '  -3.09%  '

$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").Value = '28.006.19'
$ws.Range("E21").Value = '  -2.13%  '

$ws.Range("D22").Value = "'5.308"
$ws.Range("E22").Value = '  -1.84%  '

$ws.Range("E23").Value = '  -2.93%  '

$ws.Range("D24").Value = "'2.093"

$ws.Range("D25").Value = '2.055.31'
$ws.Range("E25").Value = '  -5.90%  '

$ws.Range("D26").Value = "'153.43"
$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").Value = "'19.06"
$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("D28").Value = "'5.647"
$ws.Range("E28").Value = '  -7.53%  '

$ws.Range("D29").Value = "'1.954"
$ws.Range("E29").Value = '  -2.97%  '

$ws.Range("D30").Value = "'117.47"
$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("D31").Value = "'0.09236"
$ws.Range("E31").Value = '  -1.95%  '

$ws.Range("D32").Value = "'0.9309"
$ws.Range("E32").Value = '  -5.25%  '

$ws.Range("D33").Value = "'5.258"

$ws.Range("E34").Value = '  -2.43%  '

$ws.Range("D35").Value = "'3.293"
$ws.Range("E35").Value = '  -5.41%  '

$ws.Range("D36").Value = "'0.05833"
$ws.Range("E36").Value = '  -5.26%  '

$ws.Range("E37").Value = '  -4.00%  '

$ws.Range("D38").Value = "'1.134"
$ws.Range("E38").Value = '  -2.21%  '

$ws.Range("D39").Value = "'7.788"
$ws.Range("E39").Value = '  +2.22%  '

$ws.Range("D40").Value = "'0.5581"
$ws.Range("E40").Value = '  -2.24%  '

$ws.Range("D41").Value = "'9.850"
$ws.Range("E41").Value = '  -2.97%  '

$ws.Range("E42").Value = '  -2.28%  '

$ws.Range("E49").Value = '  -1.48%  '

$ws.Range("E50").Value = '  -0.05%  '

$ws.Range("D51").Value = "'2.323"
$ws.Range("E51").Value = '  +0.29%  '

$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").Value = "'0.07110"
$ws.Range("E43").Value = '  -0.67%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'11.57"
$ws.Range("E44").Value = '  -2.47%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.5250"
$ws.Range("E45").Value = '  -2.69%  '

$ws.Range("D46").Value = "'2.112"
$ws.Range("E46").Value = '  -10.81%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'1.824"
$ws.Range("E47").Value = '  -4.48%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = "'1.108"
$ws.Range("E48").Value = '  -11.30%  '
